$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (row 15) mirroring the structure of row 14
$row = 15

$ws.Cells.Item($row, 1).Value = 42622.891481481478
$ws.Cells.Item($row, 2).Value = 14
$ws.Cells.Item($row, 3).Value = 55
$ws.Cells.Item($row, 4).Value = 42
$ws.Cells.Item($row, 5).Value = 55
$ws.Cells.Item($row, 6).Value = 100
$ws.Cells.Item($row, 7).Value = 16793
$ws.Cells.Item($row, 8).Value = 14492
$ws.Cells.Item($row, 9).Value = 1497
$ws.Cells.Item($row, 10).Value = 264
$ws.Cells.Item($row, 11).Value = 202
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 5
$ws.Cells.Item($row, 14).Value = "Bag"

# Match the date style used by column A (same style as A2:A14) without
# introducing a brand-new style/number-format entry
$ws.Cells.Item($row - 1, 1).Copy() | Out-Null
$ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
